$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/format from the existing header cell (H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for new columns I and J (rows 2-8)
$values = @{
    2 = 9
    3 = 9
    4 = 9
    5 = 5
    6 = 7
    7 = 6
    8 = 8
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
